$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '24.647.58'
$ws.Cells.Item(2, 5).Value = '  +2.59%  '

$ws.Cells.Item(3, 4).Value = '1.697.43'
$ws.Cells.Item(3, 5).Value = '  +2.30%  '

$ws.Cells.Item(4, 4).Value = "'0.9987"
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(4, 5).Value = '  -0.89%  '

$ws.Cells.Item(5, 4).Value = "'313.94"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +1.64%  '

$ws.Cells.Item(6, 4).Value = "'0.9995"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -0.77%  '

$ws.Cells.Item(7, 4).Value = "'0.3960"
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  +1.71%  '

$ws.Cells.Item(8, 4).Value = "'0.4046"
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  +3.42%  '

$ws.Cells.Item(9, 4).Value = "'56.71"
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +16.76%  '

$ws.Cells.Item(10, 2).Value = 'Polygon'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(10, 4).Value = "'1.525"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +10.37%  '

$ws.Cells.Item(11, 2).Value = 'BinanceUSD'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(11, 4).Value = "'0.9990"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -0.83%  '

$ws.Cells.Item(12, 4).Value = "'0.08778"
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  +2.38%  '

$ws.Cells.Item(13, 4).Value = "'7.325"
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +14.22%  '

$ws.Cells.Item(14, 4).Value = "'23.18"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  +3.14%  '

$ws.Cells.Item(15, 4).Value = "'0.00001322"
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +2.73%  '

$ws.Cells.Item(16, 4).Value = "'7.612"
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  +7.44%  '

$ws.Cells.Item(17, 4).Value = '1.694.80'
$ws.Cells.Item(17, 5).Value = '  +1.90%  '

$ws.Cells.Item(18, 4).Value = "'100.63"
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  +0.21%  '

$ws.Cells.Item(19, 4).Value = "'0.07059"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  +4.57%  '

$ws.Cells.Item(20, 4).Value = "'19.53"
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +3.81%  '

$ws.Cells.Item(21, 4).Value = "'6.726"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +2.12%  '

$ws.Cells.Item(22, 4).Value = "'0.9993"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -0.77%  '

$ws.Cells.Item(23, 4).Value = "'14.18"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +4.47%  '

$ws.Cells.Item(24, 4).Value = '24.624.43'
$ws.Cells.Item(24, 5).Value = '  +2.53%  '

$ws.Cells.Item(25, 4).Value = "'3.025"
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +13.22%  '

$ws.Cells.Item(26, 4).Value = "'2.310"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -0.18%  '

$ws.Cells.Item(27, 4).Value = "'22.48"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +4.22%  '

$ws.Cells.Item(28, 4).Value = "'159.74"
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +1.49%  '

$ws.Cells.Item(29, 4).Value = "'5.193"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -0.37%  '

$ws.Cells.Item(30, 4).Value = "'133.50"
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +3.46%  '

$ws.Cells.Item(31, 4).Value = "'7.568"
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  +35.46%  '

$ws.Cells.Item(32, 4).Value = '1.877.95'
$ws.Cells.Item(32, 5).Value = '  +1.84%  '

$ws.Cells.Item(33, 5).Value = '  -2.08%  '

$ws.Cells.Item(34, 5).Value = '  +22.01%  '

$ws.Cells.Item(35, 4).Value = "'0.08567"
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  +1.26%  '

$ws.Cells.Item(36, 4).Value = "'1.963"
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +10.92%  '

$ws.Cells.Item(37, 4).Value = "'11.04"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +7.47%  '

$ws.Cells.Item(38, 4).Value = "'0.2717"
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +4.80%  '

$ws.Cells.Item(39, 4).Value = "'14.76"
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +0.00%  '

$ws.Cells.Item(40, 4).Value = "'0.02770"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +11.91%  '

$ws.Cells.Item(41, 4).Value = "'0.09042"
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  +3.37%  '

$ws.Cells.Item(42, 5).Value = '  +3.42%  '

$ws.Cells.Item(43, 4).Value = "'0.7652"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +5.07%  '

$ws.Cells.Item(44, 4).Value = "'0.7178"
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  +4.36%  '

$ws.Cells.Item(45, 4).Value = "'15.38"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +4.69%  '

$ws.Cells.Item(46, 4).Value = "'2.463"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +5.47%  '

$ws.Cells.Item(47, 4).Value = "'4.168"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +2.60%  '

$ws.Cells.Item(48, 4).Value = "'0.9989"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -0.80%  '

$ws.Cells.Item(49, 4).Value = "'1.331"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  +19.26%  '

$ws.Cells.Item(50, 4).Value = "'140.41"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +1.60%  '

$ws.Cells.Item(51, 5).Value = '  +3.03%  '
